# Entrega 3, revisora de estilo
#
# This script:
#  1) Moves the "_GoBack" bookmark from its old location (inside the
#     6th-chapter paragraph, between "...cicl" + "ic" and the rest) to the
#     very first (empty) paragraph of the document.
#  2) Accepts the proofing-tool splits (gramStart/gramEnd/spellStart/
#     spellEnd wrapped runs) by merging the surrounding runs back into a
#     single run per sentence, in five paragraphs of the "RESUMEN" section.
#
# NOTE on technique: this COM-interop runtime coalesces a freshly-edited
# run with any immediately *following* sibling run(s) that share identical
# run formatting (rPr), but never coalesces past a boundary where a format
# change was applied. We exploit that: make the whole affected run-span
# genuinely change text (placeholder trick, since an edit that nets out to
# identical text is treated as a no-op and skipped), restore the original
# text, then -- if a differently-formatted-but-actually-identical-looking
# run must stay split off from the merged run -- "flash" a formatting
# property (Bold on/off) across the remainder of the paragraph to force
# the split back, as the very last edit touching that paragraph.

$d = $word.ActiveDocument

function Merge-Run {
    param(
        [int]$ParaIndex,
        [int]$StartOffset,
        [int]$EndOffset,
        [string]$Text
    )
    $pStart = $d.Paragraphs($ParaIndex).Range.Start
    $r = $d.Range($pStart + $StartOffset, $pStart + $EndOffset)
    # Force a genuine text change (same text only flips run/proofErr
    # structure, which this engine treats as a no-op unless the literal
    # characters differ at some point), then restore the real text.
    $r.Text = $Text + "#"
    $full = $d.Paragraphs($ParaIndex).Range.Text
    $markerAt = $pStart + $StartOffset + $Text.Length
    $rMarker = $d.Range($markerAt, $markerAt + 1)
    $rMarker.Text = ""
}

function Split-After {
    param(
        [int]$ParaIndex,
        [int]$Offset
    )
    # Re-establish a run boundary at $Offset (relative to the paragraph
    # start) by toggling formatting across the remainder of the paragraph
    # and back -- this forces the host to keep that text as its own run
    # instead of re-coalescing it with the text just merged to its left.
    $p = $d.Paragraphs($ParaIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End
    $splitPoint = $pStart + $Offset
    if ($splitPoint -lt ($pEnd - 1)) {
        $rTail = $d.Range($splitPoint, $pEnd - 1)
        $rTail.Font.Bold = 1
        $rTail.Font.Bold = 0
    }
}

# --- 1. Move the _GoBack bookmark to the top of the document ---------
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()

$firstPara = $d.Paragraphs(1).Range
$d.Bookmarks.Add("_GoBack", $firstPara)

# --- 2. Paragraph 31 (RESUMEN intro) ----------------------------------
# "En el siguiente ... básica " + "de los grupo-anillos" + ", necesaria
# para el desarrollo " -> merge into one run; "de la teoría de códigos"
# (and beyond) must remain split off as a separate run.
Merge-Run -ParaIndex 31 -StartOffset 9 -EndOffset 150 `
    -Text "En el siguiente trabajo de investigación se hace un estudio detallado de la teoría básica de los grupo-anillos, necesaria para el desarrollo "
Split-After -ParaIndex 31 -Offset 150

# --- 3. Paragraph 35 (capítulo primero) -------------------------------
Merge-Run -ParaIndex 35 -StartOffset 0 -EndOffset 124 `
    -Text "El primer capítulo contiene todo el bagaje matemático que sirve de cimiento para un  estudio adecuado de los grupo-anillos. "

# --- 4. Paragraph 37 (capítulo segundo) -------------------------------
Merge-Run -ParaIndex 37 -StartOffset 0 -EndOffset 223 `
    -Text "En el segundo capítulo se da la definición de un grupo-anillo y una grupo-álgebra, caso especial del anterior. Posteriormente, se establecen las condiciones necesarias y suficientes para que un grupo-anillo sea semisimple. "

# --- 5. Paragraph 39 (capítulo tercero) --------------------------------
Merge-Run -ParaIndex 39 -StartOffset 0 -EndOffset 123 `
    -Text "En el tercer capítulo se estudia la teoría de representación de grupos y su relación con los módulos de los grupo-anillos. "

# --- 6. Paragraph 41 (capítulo cuarto) ---------------------------------
Merge-Run -ParaIndex 41 -StartOffset 0 -EndOffset 158 `
    -Text "En el cuarto capítulo se estudian algunos elementos algebraicos de un grupo-anillo como los elementos nilpotentes, los idempotentes y las unidades de torsión."

# --- 7. Paragraph 45 (capítulo sexto) -----------------------------------
# "os y mostrando ... con " + "las grupo-álgebras" + "." -> merge into one
# run; "cícl" and "ic" runs before it stay untouched.
Merge-Run -ParaIndex 45 -StartOffset 128 -EndOffset 212 `
    -Text "os y mostrando que dichos códigos tienen una fuerte conexión con las grupo-álgebras."
